$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 18.94467633333333
$ws.Range("H2").Value = 56.834029
$ws.Range("I2").Value = 0.03347881112463321
$ws.Range("J2").Value = 0.03347881112463321
$ws.Range("M2").Value = 20.29529466666667
$ws.Range("N2").Value = 60.885884
$ws.Range("O2").Value = 0.4032332285476398
$ws.Range("P2").Value = 0.4032332285476398
$ws.Range("Q2").Value = 384.4877885496263
$ws.Range("R2").Value = 3460.390096946636
$ws.Range("S2").Value = 0.01349976909772249
$ws.Range("T2").Value = 0.01349976909772249
$ws.Range("G3").Value = 18.94467633333333
$ws.Range("H3").Value = 56.834029
$ws.Range("I3").Value = 0.03347881112463321
$ws.Range("J3").Value = 0.03347881112463321
$ws.Range("O3").Value = 0.1953894087318433
$ws.Range("P3").Value = 0.1953894087318433
$ws.Range("Q3").Value = 186.3061780397145
$ws.Range("R3").Value = 1676.75560235743
$ws.Range("S3").Value = 0.00654140511068714
$ws.Range("T3").Value = 0.00654140511068714
$ws.Range("G4").Value = 18.94467633333333
$ws.Range("H4").Value = 56.834029
$ws.Range("I4").Value = 0.03347881112463321
$ws.Range("J4").Value = 0.03347881112463321
$ws.Range("M4").Value = 11.81535133333333
$ws.Range("N4").Value = 35.446054
$ws.Range("O4").Value = 0.2347510761885954
$ws.Range("P4").Value = 0.2347510761885954
$ws.Range("Q4").Value = 223.8380067746184
$ws.Range("R4").Value = 2014.542060971566
$ws.Range("S4").Value = 0.007859186941022365
$ws.Range("T4").Value = 0.007859186941022365
$ws.Range("G5").Value = 18.94467633333333
$ws.Range("H5").Value = 56.834029
$ws.Range("I5").Value = 0.03347881112463321
$ws.Range("J5").Value = 0.03347881112463321
$ws.Range("M5").Value = 8.386535
$ws.Range("N5").Value = 25.159605
$ws.Range("O5").Value = 0.1666262865319216
$ws.Range("P5").Value = 0.1666262865319216
$ws.Range("Q5").Value = 158.8801911331717
$ws.Range("R5").Value = 1429.921720198545
$ws.Range("S5").Value = 0.005578449975201218
$ws.Range("T5").Value = 0.005578449975201217
$ws.Range("I6").Value = 0.4812547190371557
$ws.Range("J6").Value = 0.4812547190371557
$ws.Range("M6").Value = 20.29529466666667
$ws.Range("N6").Value = 60.885884
$ws.Range("O6").Value = 0.4032332285476398
$ws.Range("P6").Value = 0.4032332285476398
$ws.Range("Q6").Value = 5526.975314709445
$ws.Range("R6").Value = 49742.777832385
$ws.Range("S6").Value = 0.1940578941111396
$ws.Range("T6").Value = 0.1940578941111396
$ws.Range("I7").Value = 0.4812547190371557
$ws.Range("J7").Value = 0.4812547190371557
$ws.Range("O7").Value = 0.1953894087318433
$ws.Range("P7").Value = 0.1953894087318433
$ws.Range("S7").Value = 0.09403207500207922
$ws.Range("T7").Value = 0.0940320750020792
$ws.Range("I8").Value = 0.4812547190371557
$ws.Range("J8").Value = 0.4812547190371557
$ws.Range("M8").Value = 11.81535133333333
$ws.Range("N8").Value = 35.446054
$ws.Range("O8").Value = 0.2347510761885954
$ws.Range("P8").Value = 0.2347510761885954
$ws.Range("Q8").Value = 3217.650013291389
$ws.Range("R8").Value = 28958.8501196225
$ws.Range("S8").Value = 0.1129750632148124
$ws.Range("T8").Value = 0.1129750632148124
$ws.Range("I9").Value = 0.4812547190371557
$ws.Range("J9").Value = 0.4812547190371557
$ws.Range("M9").Value = 8.386535
$ws.Range("N9").Value = 25.159605
$ws.Range("O9").Value = 0.1666262865319216
$ws.Range("P9").Value = 0.1666262865319216
$ws.Range("Q9").Value = 2283.887604602084
$ws.Range("R9").Value = 20554.98844141875
$ws.Range("S9").Value = 0.08018968670912453
$ws.Range("T9").Value = 0.08018968670912452
$ws.Range("G10").Value = 271.928284
$ws.Range("H10").Value = 815.784852
$ws.Range("I10").Value = 0.4805484928482698
$ws.Range("J10").Value = 0.4805484928482698
$ws.Range("M10").Value = 20.29529466666667
$ws.Range("N10").Value = 60.885884
$ws.Range("O10").Value = 0.4032332285476398
$ws.Range("P10").Value = 0.4032332285476398
$ws.Range("Q10").Value = 5518.864651981019
$ws.Range("R10").Value = 49669.78186782917
$ws.Range("S10").Value = 0.1937731202449102
$ws.Range("T10").Value = 0.1937731202449102
$ws.Range("G11").Value = 271.928284
$ws.Range("H11").Value = 815.784852
$ws.Range("I11").Value = 0.4805484928482698
$ws.Range("J11").Value = 0.4805484928482698
$ws.Range("O11").Value = 0.1953894087318433
$ws.Range("P11").Value = 0.1953894087318433
$ws.Range("Q11").Value = 2674.203475506094
$ws.Range("R11").Value = 24067.83127955484
$ws.Range("S11").Value = 0.09389408588460185
$ws.Range("T11").Value = 0.09389408588460185
$ws.Range("G12").Value = 271.928284
$ws.Range("H12").Value = 815.784852
$ws.Range("I12").Value = 0.4805484928482698
$ws.Range("J12").Value = 0.4805484928482698
$ws.Range("M12").Value = 11.81535133333333
$ws.Range("N12").Value = 35.446054
$ws.Range("O12").Value = 0.2347510761885954
$ws.Range("P12").Value = 0.2347510761885954
$ws.Range("Q12").Value = 3212.928212930445
$ws.Range("R12").Value = 28916.35391637401
$ws.Range("S12").Value = 0.1128092758569389
$ws.Range("T12").Value = 0.1128092758569389
$ws.Range("G13").Value = 271.928284
$ws.Range("H13").Value = 815.784852
$ws.Range("I13").Value = 0.4805484928482698
$ws.Range("J13").Value = 0.4805484928482698
$ws.Range("M13").Value = 8.386535
$ws.Range("N13").Value = 25.159605
$ws.Range("O13").Value = 0.1666262865319216
$ws.Range("P13").Value = 0.1666262865319216
$ws.Range("Q13").Value = 2280.53607125594
$ws.Range("R13").Value = 20524.82464130346
$ws.Range("S13").Value = 0.08007201086181889
$ws.Range("T13").Value = 0.08007201086181887
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.669764666666667
$ws.Range("H14").Value = 8.009294000000001
$ws.Range("I14").Value = 0.004717976989941326
$ws.Range("J14").Value = 0.004717976989941326
$ws.Range("M14").Value = 20.29529466666667
$ws.Range("N14").Value = 60.885884
$ws.Range("O14").Value = 0.4032332285476398
$ws.Range("P14").Value = 0.4032332285476398
$ws.Range("Q14").Value = 54.18366060065512
$ws.Range("R14").Value = 487.652945405896
$ws.Range("S14").Value = 0.001902445093867516
$ws.Range("T14").Value = 0.001902445093867516
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.669764666666667
$ws.Range("H15").Value = 8.009294000000001
$ws.Range("I15").Value = 0.004717976989941326
$ws.Range("J15").Value = 0.004717976989941326
$ws.Range("O15").Value = 0.1953894087318433
$ws.Range("P15").Value = 0.1953894087318433
$ws.Range("Q15").Value = 26.25506197944223
$ws.Range("R15").Value = 236.29555781498
$ws.Range("S15").Value = 0.0009218427344750774
$ws.Range("T15").Value = 0.0009218427344750774
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.669764666666667
$ws.Range("H16").Value = 8.009294000000001
$ws.Range("I16").Value = 0.004717976989941326
$ws.Range("J16").Value = 0.004717976989941326
$ws.Range("M16").Value = 11.81535133333333
$ws.Range("N16").Value = 35.446054
$ws.Range("O16").Value = 0.2347510761885954
$ws.Range("P16").Value = 0.2347510761885954
$ws.Range("Q16").Value = 31.54420751398622
$ws.Range("R16").Value = 283.897867625876
$ws.Range("S16").Value = 0.001107550175821756
$ws.Range("T16").Value = 0.001107550175821756
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.669764666666667
$ws.Range("H17").Value = 8.009294000000001
$ws.Range("I17").Value = 0.004717976989941326
$ws.Range("J17").Value = 0.004717976989941326
$ws.Range("M17").Value = 8.386535
$ws.Range("N17").Value = 25.159605
$ws.Range("O17").Value = 0.1666262865319216
$ws.Range("P17").Value = 0.1666262865319216
$ws.Range("Q17").Value = 22.39007481876334
$ws.Range("R17").Value = 201.51067336887
$ws.Range("S17").Value = 0.0007861389857769765
$ws.Range("T17").Value = 0.0007861389857769764
